# Updated validation reports to reflect current pipeline outputs
# Update Contig identifiers (column G) on Sheet1 to reflect current pipeline outputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Map of old Contig value -> new Contig value, as described by the diff.
$replacements = @{
    "Contig_2207_12.014"   = "Contig_2205_12.0141"
    "Contig_2343_11.3427"  = "Contig_2344_11.3427"
    "Contig_2000_4.63142"  = "Contig_2002_4.63142"
    "Contig_40_153.026"    = "Contig_41_153.026"
    "Contig_25_23.9329"    = "Contig_24_23.9329"
    "Contig_177_36.506"    = "Contig_176_36.506"
    "Contig_46_77.3243"    = "Contig_45_77.3243"
    "Contig_67_57.9768"    = "Contig_65_57.9768"
    "Contig_38_117.344"    = "Contig_36_117.344"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $replacements.ContainsKey([string]$val)) {
            $cell.Value2 = $replacements[[string]$val]
        }
    }
}
